# fix previous period hold-over subscription's first action being manual
# unsubscriptions not being counted
#
# Sheet "LP Rewards" (sheet4.xml) gains three new LP rows and several of
# the existing rows get updated fee/reward figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "LP Rewards"

# ---------------------------------------------------------------------
# 1) Make room for the new rows by shifting existing ones down, exactly
#    like Excel's own "Insert Sheet Rows" command would.
# ---------------------------------------------------------------------
# A brand-new LP (0x8F1c51...) now shows up right after row 5.
$ws.Rows.Item(6).Insert()

# Two more brand-new LPs (0xa614A8... and 0x93Cf0a...) show up right
# before the former rows 9/10 (now 11/12 after the first insert).
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# ---------------------------------------------------------------------
# Helper: write a value into a cell while preserving its "text" storage
# (every cell on this sheet is stored as text, even the numeric-looking
# ones) instead of letting it be auto-coerced into a Number.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------
# 2) Updated figures on rows that already existed.
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("D2") "474850.38"

Set-TextValue $ws.Range("D4") "72410.99"

Set-TextValue $ws.Range("B5") "0.003078134515545425"
Set-TextValue $ws.Range("C5") "1387.81"
Set-TextValue $ws.Range("D5") "5893.47"
Set-TextValue $ws.Range("E5") "207359"

# ---------------------------------------------------------------------
# 3) Fill in the three newly inserted rows. The lpAddress column (A) is
#    hex/hexadecimal-looking text (contains letters) so Excel already
#    stores it as text without needing the NumberFormat nudge.
# ---------------------------------------------------------------------
# Row 6 - new LP with no fees/rewards yet.
$ws.Range("A6").Value = "0x8F1c51E98Af7C7dbB24654ACB05781E0e96e008F"
Set-TextValue $ws.Range("B6") "0"
Set-TextValue $ws.Range("C6") "0"
Set-TextValue $ws.Range("D6") "0"
Set-TextValue $ws.Range("E6") "0"

# Row 10 - new LP.
$ws.Range("A10").Value = "0xa614A83132a2e7368aDa71dA9331817c33706770"
Set-TextValue $ws.Range("B10") "0.004622415337704327"
Set-TextValue $ws.Range("C10") "1756.31"
Set-TextValue $ws.Range("D10") "0"
Set-TextValue $ws.Range("E10") "0"

# Row 11 - new LP.
$ws.Range("A11").Value = "0x93Cf0a22a26895650A8AaE960Bf85a01ec6A551C"
Set-TextValue $ws.Range("B11") "0.027949695693435708"
Set-TextValue $ws.Range("C11") "21347.04"
Set-TextValue $ws.Range("D11") "89677.79"
Set-TextValue $ws.Range("E11") "3155270"

# ---------------------------------------------------------------------
# 4) The row that used to be r9 (0x3663b2...) is now r12 and its reward
#    figure also changed.
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("D12") "5315.49"
